$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: "iPhone" / "$0" (both as text, matching the existing
# shared-string-backed text cells in the table above them)
$ws.Range("A5").Value = "iPhone"

# Force B5 to be stored as literal text "$0" rather than being parsed as
# a currency value of 0 - set the cell to Text format before assigning.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "$0"
# Restore General number format so no stray per-cell style lingers beyond
# what's needed (cell ends up using the default style again).
$ws.Range("B5").NumberFormat = "General"

# Move/select the active cell to B5, matching the saved selection state.
$ws.Range("B5").Select()
